$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value  = 4216.128627752491
$ws.Range("C3").Value  = 4216.128627752491
$ws.Range("C4").Value  = 4134.683637032755
$ws.Range("C5").Value  = 4134.683637032755
$ws.Range("C6").Value  = 4087.816410299291
$ws.Range("C7").Value  = 4086.72813047241
$ws.Range("C8").Value  = 4086.72813047241
$ws.Range("C9").Value  = 4086.72813047241
$ws.Range("C10").Value = 4060.675531269419
$ws.Range("C11").Value = 3921.761139931309
$ws.Range("C12").Value = 3915.688793922054
